$d = $word.ActiveDocument

# Replace the placeholder "Kim Lam ()," entry with the student ID filled in:
# "Kim Lam ()," -> "Kim Lam (19823013),"
$d.Content.Find.Execute("Kim Lam (),", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Kim Lam (19823013),", 2)
